# Rerun model with Eigenvalues as ref
# Updates PC1 (column M) and PC2 (column N) values for each row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 2.726540799331899
$ws.Range("N2").Value = -0.6011343761358734
$ws.Range("M3").Value = 2.016045188666966
$ws.Range("N3").Value = -0.6345055137861857
$ws.Range("M4").Value = 2.482827041370223
$ws.Range("N4").Value = 0.4403078241362656
$ws.Range("M5").Value = 1.937743539785606
$ws.Range("N5").Value = -2.304290298803274
$ws.Range("M6").Value = 2.380447653508261
$ws.Range("N6").Value = -1.540343857027132
$ws.Range("M7").Value = 1.743798000363266
$ws.Range("N7").Value = 0.7509337052057601
$ws.Range("M8").Value = 1.69804703729039
$ws.Range("N8").Value = -0.5071567395949607
$ws.Range("M9").Value = 1.51113383848651
$ws.Range("N9").Value = 0.09182366035968278
$ws.Range("M10").Value = 1.766247821525466
$ws.Range("N10").Value = 1.255427573142105
$ws.Range("M11").Value = 1.832616256344844
$ws.Range("N11").Value = 1.219471586030356
$ws.Range("M12").Value = 1.109896023217602
$ws.Range("N12").Value = -0.5060339354598691
$ws.Range("M13").Value = 1.117460487216756
$ws.Range("N13").Value = -1.990620319607662
$ws.Range("M14").Value = 1.47380944131921
$ws.Range("N14").Value = -0.8287384980605181
$ws.Range("M15").Value = 0.9601669596683547
$ws.Range("N15").Value = -0.4547794531366097
$ws.Range("M16").Value = 0.8569969943998081
$ws.Range("N16").Value = -0.03661183848374457
$ws.Range("M17").Value = 0.06179640846755643
$ws.Range("N17").Value = 4.421600526028747
$ws.Range("M18").Value = 0.4477960896050738
$ws.Range("N18").Value = 2.132266818982213
$ws.Range("M19").Value = -0.03471211841238647
$ws.Range("N19").Value = -1.822425721826589
$ws.Range("M20").Value = 0.5242369474050622
$ws.Range("N20").Value = 0.9902466754491056
$ws.Range("M21").Value = -0.4803275668441702
$ws.Range("N21").Value = 0.1231207365779711
$ws.Range("M22").Value = -0.6176040000585149
$ws.Range("N22").Value = -0.01932303281527617
$ws.Range("M23").Value = -0.2281719767641495
$ws.Range("N23").Value = 0.6206050253861982
$ws.Range("M24").Value = -0.768397556496339
$ws.Range("N24").Value = -0.2251153793747109
$ws.Range("M25").Value = -0.4555226062042685
$ws.Range("N25").Value = -0.602867727774186
$ws.Range("M26").Value = -1.151055804845307
$ws.Range("N26").Value = -1.5214252792935
$ws.Range("M27").Value = -1.13268492108448
$ws.Range("N27").Value = -1.979855779964381
$ws.Range("M28").Value = -1.216094983469689
$ws.Range("N28").Value = 2.364078146885753
$ws.Range("M29").Value = -1.766923878431851
$ws.Range("N29").Value = -1.413642542745665
$ws.Range("M30").Value = -1.92785212602099
$ws.Range("N30").Value = 0.5230262448985582
$ws.Range("M31").Value = -2.355787045532871
$ws.Range("N31").Value = 2.633818040410711
$ws.Range("M32").Value = -2.824925024395875
$ws.Range("N32").Value = 0.708853052214381
$ws.Range("M33").Value = -3.749062012438849
$ws.Range("N33").Value = -1.245588494951435
$ws.Range("M34").Value = -10.46091115623349
$ws.Range("N34").Value = -0.6436418939369377
$ws.Range("M36").Value = 2.522426249260381
$ws.Range("N36").Value = 0.6025210670707041
